$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 781 (the 2026/12/29 block),
# shifting everything below it down by two rows so the sheet grows from
# A1:D822 to A1:D824.
$ws.Rows("781:782").Insert()

# Populate the two freshly inserted rows with the new 2026/02/06 entries.
# The date column is kept as plain text (like every other row in the sheet),
# so force a text number format before assigning it to stop Excel's
# automatic date-literal detection, then strip the format back off so the
# cell ends up with no explicit style - matching the rest of the sheet.
$ws.Range("A781").NumberFormat = "@"
$ws.Range("A781").Value = "2026/02/06"
$ws.Range("A781").ClearFormats()
$ws.Range("B781").Value = "金"
$ws.Range("C781").Value = 18
$ws.Range("D781").Value = 201

$ws.Range("A782").NumberFormat = "@"
$ws.Range("A782").Value = "2026/02/06"
$ws.Range("A782").ClearFormats()
$ws.Range("B782").Value = "金"
$ws.Range("C782").Value = 22
$ws.Range("D782").Value = 201
